$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.5170908230520265
$ws.Range("E2").Value = 0.5170908230520265

# Row 3
$ws.Range("D3").Value = 0.0001316036549504693
$ws.Range("E3").Value = 0.0001316036549504693

# Row 4
$ws.Range("D4").Value = [double]"8.94542481432083E-08"
$ws.Range("E4").Value = [double]"8.94542481432083E-08"

# Row 5
$ws.Range("D5").Value = 0.0119347304612219
$ws.Range("E5").Value = 0.0119347304612219

# Row 6
$ws.Range("D6").Value = 0.9672422412347126
$ws.Range("E6").Value = 0.9672422412347126

# Row 7
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 0

# Row 8
$ws.Range("D8").Value = 0.9397938006196865
$ws.Range("E8").Value = 0.06020619938031346

# Row 9
$ws.Range("D9").Value = 0.9999998271621057
$ws.Range("E9").Value = [double]"1.728378943299802E-07"

# Row 10
$ws.Range("C10").Value = $false
$ws.Range("D10").Value = 0.1711863189062275
$ws.Range("E10").Value = 0.8288136810937725

# Row 11
$ws.Range("D11").Value = 0.9216864756624517
$ws.Range("E11").Value = 0.07831352433754835
$ws.Range("F11").Value = 0.6067328453063965
$ws.Range("G11").Value = 0.7

# Row 12
$ws.Range("D12").Value = 0.9713708463916125
$ws.Range("E12").Value = 0.9713708463916125

# Row 13
$ws.Range("D13").Value = 0.001597448136458384
$ws.Range("E13").Value = 0.001597448136458384

# Row 14
$ws.Range("D14").Value = [double]"5.22854067346144E-07"
$ws.Range("E14").Value = [double]"5.22854067346144E-07"

# Row 15
$ws.Range("D15").Value = 0.01766223164804219
$ws.Range("E15").Value = 0.01766223164804219

# Row 16
$ws.Range("D16").Value = 0.4302958837927634
$ws.Range("E16").Value = 0.4302958837927634

# Row 17
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 0

# Row 18
$ws.Range("D18").Value = 0.762370498476248
$ws.Range("E18").Value = 0.237629501523752

# Row 19
$ws.Range("D19").Value = 0.9999999521487973
$ws.Range("E19").Value = [double]"4.785120266692644E-08"

# Row 20
$ws.Range("C20").Value = $false
$ws.Range("D20").Value = 0.03402114179950404
$ws.Range("E20").Value = 0.965978858200496

# Row 21
$ws.Range("D21").Value = 0.9851901291616872
$ws.Range("E21").Value = 0.01480987083831276
$ws.Range("F21").Value = 0.7802404165267944
$ws.Range("G21").Value = 0.8
